$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "1.000", "28.243.96") that must stay text.
# Temporarily force a Text number format so Excel does not auto-convert the
# assigned string into a numeric value, then clear the format again so the
# cell keeps its original (default) style - only the stored value/type changes.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.265.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4040"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.05"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08442"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.048"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.941.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.462"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.086"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06595"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.741"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.264.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.93%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.138.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.762"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -8.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.126"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9776"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09646"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.458"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.554"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.636"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.827"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02302"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.263"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06174"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6159"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1908"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.304"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5850"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.031"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.433"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06903"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.40%  "
